$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 6778.2
$ws.Range("I70").Value = 1892
$ws.Range("J70").Value = 7999.75
$ws.Range("K70").Value = 5676
$ws.Range("L70").Value = 23999.25
$ws.Range("M70").Value = -5406
$ws.Range("N70").Value = -24539.25

$ws.Range("H73").Value = 6778.2
$ws.Range("I73").Value = 1892
$ws.Range("J73").Value = 7999.75
$ws.Range("K73").Value = 5676
$ws.Range("L73").Value = 23999.25
$ws.Range("M73").Value = -4740
$ws.Range("N73").Value = -25871.25

$ws.Range("H88").Value = 1167.2222
$ws.Range("I88").Value = 1452.3334
$ws.Range("J88").Value = 1024.6666
$ws.Range("K88").Value = 1452.3334
$ws.Range("L88").Value = 1024.6666
$ws.Range("M88").Value = -1046.3334
$ws.Range("N88").Value = -1836.6666

$ws.Range("H91").Value = 1167.2222
$ws.Range("I91").Value = 1452.3334
$ws.Range("J91").Value = 1024.6666
$ws.Range("K91").Value = 1452.3334
$ws.Range("L91").Value = 1024.6666
$ws.Range("M91").Value = -48.33339999999998
$ws.Range("N91").Value = -3832.6666

$ws.Range("H101").Value = 491.375
$ws.Range("J101").Value = 431.66666
$ws.Range("L101").Value = 1294.99998
$ws.Range("N101").Value = -4538.999980000001

$ws.Range("H111").Value = 1325.3334
$ws.Range("I111").Value = 1238
$ws.Range("K111").Value = 3714
$ws.Range("M111").Value = -647

$ws.Range("H125").Value = 2214.25
$ws.Range("I125").Value = 2069.389
$ws.Range("J125").Value = 2648.8333
$ws.Range("K125").Value = 18624.501
$ws.Range("L125").Value = 23839.4997
$ws.Range("M125").Value = -16164.501
$ws.Range("N125").Value = -28759.4997

$ws.Range("H132").Value = 14681.667
$ws.Range("I132").Value = 14138
$ws.Range("K132").Value = 42414
$ws.Range("M132").Value = -39884

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5546.722
$ws.Range("I2").Value = 3651.2144
$ws.Range("K2").Value = 3651.2144
$ws.Range("M2").Value = -3538.2144

$ws.Range("H110").Value = 3138.3572
$ws.Range("I110").Value = 3528
$ws.Range("K110").Value = 3528
$ws.Range("M110").Value = -1483

$ws.Range("H116").Value = 5546.722
$ws.Range("I116").Value = 3651.2144
$ws.Range("K116").Value = 3651.2144
$ws.Range("M116").Value = -1357.2144

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5546.722
$ws.Range("I3").Value = 3651.2144
$ws.Range("K3").Value = 3651.2144
$ws.Range("M3").Value = -3537.2144

$ws.Range("H86").Value = 3708.353
$ws.Range("I86").Value = 1572.125
$ws.Range("J86").Value = 5607.222
$ws.Range("K86").Value = 1572.125
$ws.Range("L86").Value = 5607.222
$ws.Range("M86").Value = -449.125
$ws.Range("N86").Value = -7853.222

$ws.Range("H89").Value = 3708.353
$ws.Range("I89").Value = 1572.125
$ws.Range("J89").Value = 5607.222
$ws.Range("K89").Value = 7860.625
$ws.Range("L89").Value = 28036.11
$ws.Range("M89").Value = -2244.625
$ws.Range("N89").Value = -39268.11

$ws.Range("H94").Value = 898.75
$ws.Range("I94").Value = 875
$ws.Range("J94").Value = 922.5
$ws.Range("K94").Value = 875
$ws.Range("L94").Value = 922.5
$ws.Range("M94").Value = -424
$ws.Range("N94").Value = -1824.5

$ws.Range("H107").Value = 3904.76
$ws.Range("I107").Value = 1108
$ws.Range("K107").Value = 1108
$ws.Range("M107").Value = 812

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()

$ws.Range("H51").Value = 49978
$ws.Range("J51").Value = 49978
$ws.Range("L51").Value = 49978
$ws.Range("N51").Value = -51450

$ws.Range("H61").Value = 49978
$ws.Range("J61").Value = 49978
$ws.Range("L61").Value = 49978
$ws.Range("N61").Value = -50674

$ws.Range("H80").Value = 94036.25
$ws.Range("J80").Value = 94036.25
$ws.Range("L80").Value = 94036.25
$ws.Range("N80").Value = -96282.25

$ws.Range("H83").Value = 94036.25
$ws.Range("J83").Value = 94036.25
$ws.Range("L83").Value = 282108.75
$ws.Range("N83").Value = -293340.75

$ws.Range("H107").Value = 276
$ws.Range("I107").Value = 282.42856
$ws.Range("J107").Value = 246
$ws.Range("K107").Value = 282.42856
$ws.Range("L107").Value = 246
$ws.Range("M107").Value = 1637.57144
$ws.Range("N107").Value = -4086

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value = 1379.5
$ws.Range("I122").Value = 1249
$ws.Range("K122").Value = 3747
$ws.Range("M122").Value = -1297

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("M119").ClearContents()

$ws.Range("H132").Value = 1887.25
$ws.Range("J132").Value = 1887.25
$ws.Range("L132").Value = 16985.25
$ws.Range("N132").Value = -22045.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1107.9166
$ws.Range("I102").Value = 1107.9166
$ws.Range("K102").Value = 1107.9166
$ws.Range("M102").Value = 514.0834

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2705.3157
$ws.Range("I61").Value = 1293.4
$ws.Range("K61").Value = 1293.4
$ws.Range("M61").Value = -1091.4

$ws.Range("H68").Value = 8514.177
$ws.Range("I68").Value = 8249.666999999999
$ws.Range("J68").Value = 8570.857
$ws.Range("K68").Value = 8249.666999999999
$ws.Range("L68").Value = 8570.857
$ws.Range("M68").Value = -7500.666999999999
$ws.Range("N68").Value = -10068.857

$ws.Range("H69").Value = 34499
$ws.Range("J69").Value = 34499
$ws.Range("L69").Value = 34499
$ws.Range("N69").Value = -36121

$ws.Range("H71").Value = 8514.177
$ws.Range("I71").Value = 8249.666999999999
$ws.Range("J71").Value = 8570.857
$ws.Range("K71").Value = 41248.335
$ws.Range("L71").Value = 42854.285
$ws.Range("M71").Value = -37504.335
$ws.Range("N71").Value = -50342.285

$ws.Range("H72").Value = 34499
$ws.Range("J72").Value = 34499
$ws.Range("L72").Value = 103497
$ws.Range("N72").Value = -111609

$ws.Range("H88").Value = 85000
$ws.Range("J88").Value = 85000
$ws.Range("L88").Value = 85000
$ws.Range("N88").Value = -85856

$ws.Range("H91").Value = 85000
$ws.Range("J91").Value = 85000
$ws.Range("L91").Value = 85000
$ws.Range("N91").Value = -87964

$ws.Range("H113").Value = 2705.3157
$ws.Range("I113").Value = 1293.4
$ws.Range("K113").Value = 1293.4
$ws.Range("M113").Value = 876.5999999999999

$ws.Range("H136").Value = 3488.8
$ws.Range("I136").Value = 2871.25
$ws.Range("K136").Value = 8613.75
$ws.Range("M136").Value = -6063.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7448.9
$ws.Range("I62").Value = 3349.75
$ws.Range("J62").Value = 10181.667
$ws.Range("K62").Value = 3349.75
$ws.Range("L62").Value = 10181.667
$ws.Range("M62").Value = -2725.75
$ws.Range("N62").Value = -11429.667

$ws.Range("H65").Value = 7448.9
$ws.Range("I65").Value = 3349.75
$ws.Range("J65").Value = 10181.667
$ws.Range("K65").Value = 16748.75
$ws.Range("L65").Value = 50908.335
$ws.Range("M65").Value = -13628.75
$ws.Range("N65").Value = -57148.335

$ws.Range("H69").Value = 22635
$ws.Range("J69").Value = 22635
$ws.Range("L69").Value = 22635
$ws.Range("N69").Value = -24133

$ws.Range("H72").Value = 22635
$ws.Range("J72").Value = 22635
$ws.Range("L72").Value = 67905
$ws.Range("N72").Value = -75393

$ws.Range("H122").Value = 2079.8333
$ws.Range("I122").Value = 2079.8333
$ws.Range("K122").Value = 6239.499899999999
$ws.Range("M122").Value = -3789.499899999999
